$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New review rows (11-14) mirror the layout/format of row 10 (the last
# populated review row). Copy that row's formatting+values down into each
# new row first, then overwrite the cell values that differ per row.
$ws.Range("A10:G10").Copy()
$ws.Range("A11:G11").PasteSpecial(-4104)
$ws.Range("A10:G10").Copy()
$ws.Range("A12:G12").PasteSpecial(-4104)
$ws.Range("A10:G10").Copy()
$ws.Range("A13:G13").PasteSpecial(-4104)
$ws.Range("A10:G10").Copy()
$ws.Range("A14:G14").PasteSpecial(-4104)

# Column A's cell style (bold "appid" look) doesn't survive the row copy
# above, so reassert its font explicitly for each new row.
$ws.Range("A11").Font.Name = "Mangal"
$ws.Range("A11").Font.Size = 10
$ws.Range("A12").Font.Name = "Mangal"
$ws.Range("A12").Font.Size = 10
$ws.Range("A13").Font.Name = "Mangal"
$ws.Range("A13").Font.Size = 10
$ws.Range("A14").Font.Name = "Mangal"
$ws.Range("A14").Font.Size = 10

# --- Row 11 ---
$ws.Range("C11").Value = "nevilgreen12@gmail.com"
$ws.Range("D11").Value = "vikicrestina@gmail.com"
$ws.Range("F11").Value = "tremendous effect in my life after this guide. Thank you bitcoin"

# --- Row 12 ---
$ws.Range("C12").Value = "snizzvered@gmail.com"
$ws.Range("D12").Value = "krigelron@gmail.com"
$ws.Range("F12").Value = "gold info, good"

# --- Row 13 ---
$ws.Range("C13").Value = "redvelvetmichael@gmail.com"
$ws.Range("D13").Value = "veredsnir12@gmail.com"
$ws.Range("F13").Value = "good app"

# --- Row 14 ---
$ws.Range("C14").Value = "veredsnir12@gmail.com"
$ws.Range("D14").Value = "kevinkors122@gmail.com"
$ws.Range("F14").Value = "informative and special for beginners"

# Hyperlinks for the email columns (C/D) of the new rows.
$ws.Hyperlinks.Add($ws.Range("C11"), "mailto:nevilgreen12@gmail.com", "", "", "nevilgreen12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D11"), "mailto:vikicrestina@gmail.com", "", "", "vikicrestina@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C12"), "mailto:snizzvered@gmail.com", "", "", "snizzvered@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D12"), "mailto:krigelron@gmail.com", "", "", "krigelron@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C13"), "mailto:redvelvetmichael@gmail.com", "", "", "redvelvetmichael@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D13"), "mailto:veredsnir12@gmail.com", "", "", "veredsnir12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C14"), "mailto:veredsnir12@gmail.com", "", "", "veredsnir12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D14"), "mailto:kevinkors122@gmail.com", "", "", "kevinkors122@gmail.com")

# Adding a hyperlink restyles the cell with the hyperlink theme; restore
# each cell's original (row-10-matching) formatting afterwards.
$ws.Range("C10:D10").Copy()
$ws.Range("C11:D11").PasteSpecial(-4122)
$ws.Range("C10:D10").Copy()
$ws.Range("C12:D12").PasteSpecial(-4122)
$ws.Range("C10:D10").Copy()
$ws.Range("C13:D13").PasteSpecial(-4122)
$ws.Range("C10:D10").Copy()
$ws.Range("C14:D14").PasteSpecial(-4122)

# Restore the active selection to the cell that was selected when the
# workbook was last saved.
$ws.Range("C14").Select()
